$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkTypes")

$ws.Range("G3").Value = "кол-во факт"
$ws.Range("H3").Value = "сумма факт"

$ws.Range("G9").Value = "кол-во факт"
$ws.Range("H9").Value = "сумма факт"

$ws.Range("B15").Value = "ИТОГО по объекту"
